$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns A, B, D, E to be stored as literal text (matches source inlineStr cells)
# by setting a Text number format before assignment, then clearing the format afterward
# so no extra style survives in the saved file (matches the unstyled target cells).
$ws.Range("A73:B89").NumberFormat = "@"
$ws.Range("D73:E89").NumberFormat = "@"

$ws.Range("A73").Value = '2021-04-06'
$ws.Range("B73").Value = '03:08:44'
$ws.Range("C73").Value = 2
$ws.Range("D73").Value = '3430.0'
$ws.Range("E73").Value = '3100.0'
$ws.Range("H73").Value = 'InService'

$ws.Range("A74").Value = '2021-04-06'
$ws.Range("B74").Value = '03:10:39'
$ws.Range("C74").Value = 2
$ws.Range("D74").Value = '5255.0'
$ws.Range("E74").Value = '5626.0'
$ws.Range("H74").Value = 'InService'

$ws.Range("A75").Value = '2021-04-06'
$ws.Range("B75").Value = '03:11:15'
$ws.Range("C75").Value = 1.8333
$ws.Range("D75").Value = '3486.0'
$ws.Range("E75").Value = '3028.0'
$ws.Range("H75").Value = 'InService'

$ws.Range("A76").Value = '2021-04-06'
$ws.Range("B76").Value = '03:12:53'
$ws.Range("C76").Value = 2.1667
$ws.Range("D76").Value = '3510.0'
$ws.Range("E76").Value = '3204.0'
$ws.Range("H76").Value = 'InService'

$ws.Range("A77").Value = '2021-04-06'
$ws.Range("B77").Value = '03:13:19'
$ws.Range("C77").Value = 1.8333
$ws.Range("D77").Value = '3274.0'
$ws.Range("E77").Value = '2924.0'
$ws.Range("H77").Value = 'InService'

$ws.Range("A78").Value = '2021-04-06'
$ws.Range("B78").Value = '03:13:39'
$ws.Range("C78").Value = 1.8333
$ws.Range("D78").Value = '3274.0'
$ws.Range("E78").Value = '2924.0'
$ws.Range("H78").Value = 'InService'

$ws.Range("A79").Value = '2021-04-06'
$ws.Range("B79").Value = '03:14:00'
$ws.Range("C79").Value = 2.1667
$ws.Range("D79").Value = '6230.0'
$ws.Range("E79").Value = '11206.0'
$ws.Range("H79").Value = 'InService'

$ws.Range("A80").Value = '2021-04-06'
$ws.Range("B80").Value = '03:14:21'
$ws.Range("C80").Value = 2.1667
$ws.Range("D80").Value = '6230.0'
$ws.Range("E80").Value = '11206.0'
$ws.Range("H80").Value = 'InService'

$ws.Range("A81").Value = '2021-04-06'
$ws.Range("B81").Value = '03:16:35'
$ws.Range("C81").Value = 2.0339
$ws.Range("D81").Value = '3330.0'
$ws.Range("E81").Value = '4436.0'
$ws.Range("H81").Value = 'InService'

$ws.Range("A82").Value = '2021-04-06'
$ws.Range("B82").Value = '03:16:56'
$ws.Range("C82").Value = 2.0339
$ws.Range("D82").Value = '3330.0'
$ws.Range("E82").Value = '4436.0'
$ws.Range("H82").Value = 'InService'

$ws.Range("A83").Value = '2021-04-06'
$ws.Range("B83").Value = '03:18:48'
$ws.Range("C83").Value = 2
$ws.Range("D83").Value = '3280.0'
$ws.Range("E83").Value = '2964.0'
$ws.Range("H83").Value = 'InService'

$ws.Range("A84").Value = '2021-04-06'
$ws.Range("B84").Value = '03:19:28'
$ws.Range("C84").Value = 1.8033
$ws.Range("D84").Value = '5392.0'
$ws.Range("E84").Value = '5631.0'
$ws.Range("H84").Value = 'InService'

$ws.Range("A85").Value = '2021-04-06'
$ws.Range("B85").Value = '03:19:49'
$ws.Range("C85").Value = 1.8033
$ws.Range("D85").Value = '5392.0'
$ws.Range("E85").Value = '5631.0'
$ws.Range("H85").Value = 'InService'

$ws.Range("A86").Value = '2021-04-06'
$ws.Range("B86").Value = '03:20:10'
$ws.Range("C86").Value = 2.0339
$ws.Range("D86").Value = '3340.0'
$ws.Range("E86").Value = '3024.0'
$ws.Range("H86").Value = 'InService'

$ws.Range("A87").Value = '2021-04-06'
$ws.Range("B87").Value = '03:20:31'
$ws.Range("C87").Value = 2.0339
$ws.Range("D87").Value = '3340.0'
$ws.Range("E87").Value = '3024.0'
$ws.Range("H87").Value = 'InService'

$ws.Range("A88").Value = '2021-04-06'
$ws.Range("B88").Value = '03:20:51'
$ws.Range("C88").Value = 2.0339
$ws.Range("D88").Value = '3340.0'
$ws.Range("E88").Value = '3024.0'
$ws.Range("H88").Value = 'InService'

$ws.Range("A89").Value = '2021-04-06'
$ws.Range("B89").Value = '03:21:12'
$ws.Range("C89").Value = 2
$ws.Range("D89").Value = '3364.0'
$ws.Range("E89").Value = '3056.0'
$ws.Range("H89").Value = 'InService'

$ws.Range("A73:B89").ClearFormats()
$ws.Range("D73:E89").ClearFormats()
